# Applies two spelling/merge corrections:
#   Slide 4: "Python " + "Scikit"(err) + "-Learn"  ->  single run "Python Scikit-Learn"
#   Slide 5: "Week Type (weekday, " + "Saurday"(err) + ", or Sunday)"
#            -> single run "Week Type (weekday, Saturday, or Sunday)" (fixes the typo too)
#
# The original runs are merged into one run by: wiping out the middle
# (misspelled / err-flagged) run, rewriting the first run's text to the
# full corrected string, and then wiping out the now-trailing leftover
# run. Setting a run's .Text to "" removes that run element entirely,
# which is what lets the three runs collapse back down to one - matching
# the target OOXML (single <a:r> with the original first run's <a:rPr>).

$p = $ppt.ActivePresentation

function Merge-Paragraph-Runs {
    param($TextRange, $ParaIndex, $NewText)

    $para = $TextRange.Paragraphs($ParaIndex, 1)

    # Remove the second (middle) run.
    $run2 = $para.Runs(2, 1)
    $run2.Text = ""

    # Put the full corrected text into what is now the first run -
    # this preserves its original rPr (lang/sz/kern/dirty) and drops
    # any err="1" flag that belonged to the old middle run.
    $para = $TextRange.Paragraphs($ParaIndex, 1)
    $run1 = $para.Runs(1, 1)
    $run1.Text = $NewText

    # Remove the leftover tail run (the old third run's text is now
    # duplicated at the end of run1).
    $para = $TextRange.Paragraphs($ParaIndex, 1)
    $run2b = $para.Runs(2, 1)
    $run2b.Text = ""
}

# --- Slide 4: "Python Scikit-Learn" bullet ---
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(3)
$tr4 = $shape4.TextFrame.TextRange
Merge-Paragraph-Runs $tr4 7 "Python Scikit-Learn"

# --- Slide 5: "Week Type (weekday, Saturday, or Sunday)" bullet ---
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(3)
$tr5 = $shape5.TextFrame.TextRange
Merge-Paragraph-Runs $tr5 11 "Week Type (weekday, Saturday, or Sunday)"
